# Update row 6 (J6:AS6) on "strategy_id-0" from 1 to 0.5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")
$ws.Range("J6:AS6").Value = 0.5

# Delete the three extra templated strategy sheets
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("strategy_id-5004").Delete()
$wb.Worksheets.Item("strategy_id-5007").Delete()
$wb.Worksheets.Item("strategy_id-5009").Delete()
